$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.205.80'
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.355.12'
$ws.Range("E3").Value = '  +2.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.61'
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.58'
$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.505'
$ws.Range("E7").Value = '  -1.06%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  +0.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.14'
$ws.Range("E10").Value = '  -1.73%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("E11").Value = '  +0.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.58'
$ws.Range("E12").Value = '  -3.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.122'
$ws.Range("E13").Value = '  +3.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("E14").Value = '  -1.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.716.99'
$ws.Range("E15").Value = '  +2.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.345.47'
$ws.Range("E16").Value = '  +1.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.798'
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.168.95'
$ws.Range("E18").Value = '  +1.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.20'
$ws.Range("E19").Value = '  -0.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.24'
$ws.Range("E20").Value = '  +3.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0890'
$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.19'
$ws.Range("E22").Value = '  +0.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.30'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  -1.12%  '

$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("E26").Value = '  +1.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.60'
$ws.Range("E27").Value = '  +0.04%  '

$ws.Range("E28").Value = '  +14.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.18'
$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.26'
$ws.Range("E30").Value = '  -2.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.02'
$ws.Range("E32").Value = '  +0.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0724'
$ws.Range("E33").Value = '  +3.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.20'
$ws.Range("E34").Value = '  -2.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.85'
$ws.Range("E35").Value = '  +5.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.36'
$ws.Range("E36").Value = '  -1.42%  '

$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("E38").Value = '  +0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.44'
$ws.Range("E39").Value = '  +13.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.76'
$ws.Range("E40").Value = '  +1.86%  '

$ws.Range("E41").Value = '  -0.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.73'
$ws.Range("E42").Value = '  -32.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.941.78'
$ws.Range("E43").Value = '  -1.00%  '

$ws.Range("E44").Value = '  +0.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.12'
$ws.Range("E45").Value = '  +3.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.59'
$ws.Range("E46").Value = '  -8.46%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.580.66'
$ws.Range("E48").Value = '  +2.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.84'
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("E50").Value = '  -3.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.29'
$ws.Range("E51").Value = '  +0.94%  '
